# SRS workbook update:
#  - Enhancements: new row 28 (check-in/booking select option task)
#  - queries: several new Q&A rows about mongo aggregate queries + web service note
#  - view/selection housekeeping: make "queries" the active sheet/tab

$wb = $excel.ActiveWorkbook
$wsEnh = $wb.Worksheets.Item("Enhancements")
$wsQ = $wb.Worksheets.Item("queries")

# ---------------------------------------------------------------------------
# 1. Enhancements: append row 28, same look & feel as row 27 (copy formats)
# ---------------------------------------------------------------------------
$wsEnh.Range("A27:F27").Copy()
$wsEnh.Range("A28:F28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsEnh.Range("A28").Value = 42635
$wsEnh.Range("B28").Value = "select option for check in and booking"
$wsEnh.Range("C28").Value = "select option in the check in page so that for each room the user can make check in (in case of same date) or booking (for future) in the check in tab`nNote that this value the value must be either CHECKED-IN or BOOKED as per data coming from server and must be send with key bookingStatus as mentioned in json of /transaction(post)"
$wsEnh.Range("D28").Value = "rajashree"
$wsEnh.Range("E28").Value = 42634
$wsEnh.Range("F28").Value = 42634
$wsEnh.Rows.Item(28).RowHeight = 150

# ---------------------------------------------------------------------------
# 2. queries: restyle header/first row + append new question/answer rows
# ---------------------------------------------------------------------------

# Header row + row 2 pick up the same cell style already used elsewhere on
# the sheet (header style / wrap-text body style) instead of their old
# one-off styles.
$wsQ.Range("A1:C1").WrapText = $true
$wsQ.Range("A2").WrapText = $true
$wsQ.Range("B2:C2").WrapText = $true

$wsQ.Range("B3").Value = "get the price matching with the array of the rooms"
$wsQ.Range("C3").Value = "db.getCollection('tranctions').aggregate([{ "

$wsQ.Range("C4").Value = "        `$match : {"
$wsQ.Range("C5").Value = "            'roomsDetails.price':100"
$wsQ.Range("C6").Value = "            "
$wsQ.Range("C7").Value = "         }"
$wsQ.Range("C8").Value = "    }])"

$wsQ.Range("A10").Value = 42635
$wsQ.Range("B10").Value = "Note if we want to get data matching key value from the array of objects then we have to use aggregrate"
$wsQ.Range("B11").Value = "if we want to update data in array of objects then we have to use `$"

$wsQ.Range("A13").Value = 42635
$wsQ.Range("B13").Value = "in case we want to use aggregrate and populate then first we have to get the data from aggregrate and then populate it"

$wsQ.Range("A15").Value = 42635
$wsQ.Range("B15").Value = "updated the web service for the get transaction where the parameres are remain the same but the response is different"

# New rows reuse the plain wrap-text body style (same as B2/C2 above) and
# the date style used for A2/A10/A13/A15.
$wsQ.Range("B2:C2").Copy()
$wsQ.Range("B3:C3").PasteSpecial(-4122)
$wsQ.Range("C4").PasteSpecial(-4122)
$wsQ.Range("C5").PasteSpecial(-4122)
$wsQ.Range("C6").PasteSpecial(-4122)
$wsQ.Range("C7").PasteSpecial(-4122)
$wsQ.Range("C8").PasteSpecial(-4122)
$wsQ.Range("B10").PasteSpecial(-4122)
$wsQ.Range("B11").PasteSpecial(-4122)
$wsQ.Range("B13").PasteSpecial(-4122)
$wsQ.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ.Range("A2").Copy()
$wsQ.Range("A10").PasteSpecial(-4122)
$wsQ.Range("A13").PasteSpecial(-4122)
$wsQ.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ.Rows.Item(10).RowHeight = 30
$wsQ.Rows.Item(13).RowHeight = 30
$wsQ.Rows.Item(15).RowHeight = 30

# ---------------------------------------------------------------------------
# 3. View / selection housekeeping
# ---------------------------------------------------------------------------
$wsEnh.Range("C28").Select()
$wsQ.Range("C15").Select()
$wsQ.Activate()

Write-Output "done"
